$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row at position 3 (old row3 -> row4, old row4 -> row5)
$ws.Rows.Item(3).Insert()

# 2) Seed the new row 3 with the same "template" data as row 2 (copy/paste),
#    then fix up the row index + the two cells that differ for this test case.
$ws.Range("A2:AA2").Copy()
$ws.Range("A3:AA3").PasteSpecial()
$ws.Range("A3").Value2 = 2

# 3) Replay the sequence of submission-order-number / note edits in the same
#    order the test run produced them, so the shared-string table grows the
#    same way. Cell B2 is repeatedly updated with the latest "Submission
#    completed..." order number; E3/B3 get their one-off values in between.
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459290"
$ws.Range("E3").Value2  = "abc1"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459320"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459324"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459347"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459350"
$ws.Range("B3").Value2  = "Wrong product added for Submission"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459413"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459420"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459421"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459423"
$ws.Range("B2").Value2  = "Submission completed with order number: Request details    /    PR459450"

# 4) Match the state of row 5 (old row 4): state abbreviation now "NC".
$ws.Range("Z5").Value2 = "NC"

# 5) Update the view: scroll back to A1 and select A2 (was topLeftCell F1 / AA2).
$ws.Range("A1").Select()
$ws.Range("A2").Select()
